{"js": "const body = context.document.body;\n\n// 1. Replace the inline image in the first paragraph with the question text.\nconst pics = body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nif (pics.items.length > 0) {\n  const pic = pics.items[0];\n  const picParagraph = pic.paragraph;\n  pic.delete();\n  picParagraph.insertText(\n    \"35. Ausha placed a square plece of red filter in front of a torch as shown in diagram 1.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2. Small text/typo fixes throughout the body.\nconst replacements = [\n  [\n    \"When the torch was switched on, she noted a Square. patch of light on the\",\n    \"When the torch was switched on, she noted.a Square patch of light on the\",\n  ],\n  [\n    \"Allsha then placed a can of drink behind a red filter as shown in diagram 2.\",\n    \"Alisha then placed a can of drink behind a red filter as shown in diagram 2.\",\n  ],\n  [\". me ue. black shadow\", \". me oe, black shadow\"],\n  [\n    \"3: ped filter . Can ofdrink. . ot tignt\",\n    \"sped \\u2018filter . Can of drink. . ot tignt\",\n  ],\n  [\n    \"(b). .. Diagram 3. shows what. she saw on the. screen when the forch was\",\n    \"(b). . Diagram 3. shows: what.she saw on the screen when the torch was\",\n  ],\n  [\"Ce : [2]\", \"oo ae : [2]\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2, $null, $null, $null, $null)\n}\n\n# 1. Replace the image in paragraph 1 with text\n$shp = $d.InlineShapes(1)\n$p1 = $d.Paragraphs(1)\n$shp.Delete()\n$p1.Range.Text = \"35. Ausha placed a square plece of red filter in front of a torch as shown in diagram 1.\"\n\n# 2. Text replacements\nReplace-Text \"When the torch was switched on, she noted a Square. patch of light on the\" \"When the torch was switched on, she noted.a Square patch of light on the\"\nReplace-Text \"Allsha then placed a can of drink behind a red filter as shown in diagram 2.\" \"Alisha then placed a can of drink behind a red filter as shown in diagram 2.\"\nReplace-Text \". me ue. black shadow\" \". me oe, black shadow\"\nReplace-Text \"3: ped filter . Can ofdrink. . ot tignt\" \"sped \u2018filter . Can of drink. . ot tignt\"\nReplace-Text \"(b). .. Diagram 3. shows what. she saw on the. screen when the forch was\" \"(b). . Diagram 3. shows: what.she saw on the screen when the torch was\"\nReplace-Text \"Ce : [2]\" \"oo ae : [2]\"\n\nWrite-Output \"Done\"\n"}
